$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# Row 2
$wsForecast.Range("D2").Value = 2
$wsForecast.Range("H2").Value = 4.76
$wsForecast.Range("L2").Value = 0.96

# Row 3
$wsForecast.Range("H3").Value = 3.46
$wsForecast.Range("L3").Value = 1.2

# Row 4
$wsForecast.Range("H4").Value = 2.46
$wsForecast.Range("L4").Value = 0.9

# Row 5
$wsForecast.Range("H5").Value = 1.46
$wsForecast.Range("L5").Value = 0.87

# Row 6
$wsForecast.Range("H6").Value = 0.46
$wsForecast.Range("I6").Value = "High"
$wsForecast.Range("L6").Value = 1.16

# Row 7
$wsForecast.Range("L7").Value = 1.18

# Row 8
$wsForecast.Range("L8").Value = 0.83

# Row 9
$wsForecast.Range("L9").Value = 0.84

# Row 10
$wsForecast.Range("L10").Value = 0.9399999999999999

# Row 11
$wsForecast.Range("L11").Value = 0.9399999999999999

# Row 12
$wsForecast.Range("L12").Value = 0.99

# Row 13
$wsForecast.Range("L13").Value = 0.83

# Row 14
$wsForecast.Range("L14").Value = 0.83

# Row 15
$wsForecast.Range("D15").Value = 1
$wsForecast.Range("L15").Value = 1.19

# Row 16
$wsForecast.Range("L16").Value = 1.02

# Row 17
$wsForecast.Range("L17").Value = 0.83

# Summary sheet
$wsSummary.Range("B10").Value = "19"
